$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "RPA수행결과" column (O) to the result sheet, mirroring the
# existing "직번"(A) header/body formatting (merged header cell across the
# two header rows, plain bordered body cells for rows 3-5).

# 1) Merge the new header cell first (O1:O2), matching A1:A2 / B1:B2 / etc.
[void]$ws.Range("O1:O2").Merge()

# 2) Copy formatting from the matching cells in column A so the new column
#    reuses the existing styles instead of minting unrelated ones.
$ws.Range("A1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("O2").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("O3:O5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Set the header text for the new column.
$ws.Range("O1").Value = "RPA수행결과"

# 4) Leave the selection on the newly added header cell, like Excel would
#    after a user finishes editing it.
[void]$ws.Range("O1:O2").Select()
